$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) / Volume(1h) (E) columns with the latest crypto snapshot.
# Numeric-looking price strings are written with a leading apostrophe so Excel
# keeps them as text (matching the existing text-formatted Price column)
# instead of auto-converting them to numbers.

$ws.Range("D2").Value = "26.342.61"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "1.667.29"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.93%  "
$ws.Range("D5").Value = "'219.46"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "'0.5358"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").Value = "'0.06415"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'0.07852"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "1.674.68"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "1.894.61"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'0.5546"
$ws.Range("D16").Value = "0.0₅8200"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "'65.84"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "26.361.86"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "'4.695"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'193.63"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'10.31"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").Value = "'6.053"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").Value = "'146.39"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").Value = "'0.1233"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'7.229"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "'1.499"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("D30").Value = "'0.05876"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "'1.287"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'3.640"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "'3.286"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").Value = "'0.9724"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").Value = "'2.823"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").Value = "'2.421"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "'0.5842"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'0.8693"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("D41").Value = "'5.853"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "1.061.96"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").Value = "'105.12"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D45").Value = "1.805.99"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'57.93"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").Value = "'8.022"
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "'0.05168"
$ws.Range("E51").Value = "  +0.48%  "

# Rows 47-48: Frax and BabyDogeCoin swapped ranking positions this refresh,
# so update Coin/Link/Price/Volume for both rows accordingly.
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.016"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -6.97%  "
